$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 678.2308
$ws.Range("I41").Value = 477.75
$ws.Range("K41").Value = 477.75
$ws.Range("M41").Value = -37.75

$ws.Range("H112").Value = 2107.625
$ws.Range("J112").Value = 3433.5
$ws.Range("L112").Value = 10300.5
$ws.Range("N112").Value = -12516.5

$ws.Range("H137").Value = 4072.975
$ws.Range("I137").Value = 5166.9644
$ws.Range("J137").Value = 1520.3334
$ws.Range("K137").Value = 15500.8932
$ws.Range("L137").Value = 4561.0002
$ws.Range("M137").Value = -12950.8932
$ws.Range("N137").Value = -9661.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1715.26
$ws.Range("I32").Value = 1639.4949
$ws.Range("K32").Value = 1639.4949
$ws.Range("M32").Value = -1352.4949

$ws.Range("H61").Value = 5184.061
$ws.Range("I61").Value = 4561.515
$ws.Range("J61").Value = 6468.0625
$ws.Range("K61").Value = 4561.515
$ws.Range("L61").Value = 6468.0625
$ws.Range("M61").Value = -4349.515
$ws.Range("N61").Value = -6892.0625

$ws.Range("H88").Value = 1544.3334
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 1599.75
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1599.75
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2411.75

$ws.Range("H91").Value = 1544.3334
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 1599.75
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1599.75
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4407.75

$ws.Range("H136").Value = 5184.061
$ws.Range("I136").Value = 4561.515
$ws.Range("J136").Value = 6468.0625
$ws.Range("K136").Value = 13684.545
$ws.Range("L136").Value = 19404.1875
$ws.Range("M136").Value = -11134.545
$ws.Range("N136").Value = -24504.1875

$ws.Range("H141").Value = 147484.5
$ws.Range("J141").Value = 147484.5
$ws.Range("L141").Value = 147484.5
$ws.Range("N141").Value = -157844.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5257.2
$ws.Range("I20").Value = 5007.5
$ws.Range("K20").Value = 5007.5
$ws.Range("M20").Value = -4760.5

$ws.Range("H86").Value = 4305.6665
$ws.Range("I86").Value = 2969.5
$ws.Range("J86").Value = 4973.75
$ws.Range("K86").Value = 2969.5
$ws.Range("L86").Value = 4973.75
$ws.Range("M86").Value = -1846.5
$ws.Range("N86").Value = -7219.75

$ws.Range("H89").Value = 4305.6665
$ws.Range("I89").Value = 2969.5
$ws.Range("J89").Value = 4973.75
$ws.Range("K89").Value = 14847.5
$ws.Range("L89").Value = 24868.75
$ws.Range("M89").Value = -9231.5
$ws.Range("N89").Value = -36100.75

$ws.Range("H94").Value = 1186.6666
$ws.Range("I94").Value = 713.7241
$ws.Range("K94").Value = 713.7241
$ws.Range("M94").Value = -262.7241

$ws.Range("H99").Value = 1600.9166
$ws.Range("I99").Value = 1496.75
$ws.Range("K99").Value = 1496.75
$ws.Range("M99").Value = 1.25

$ws.Range("H105").Value = 2531.75
$ws.Range("I105").Value = 2504.5386
$ws.Range("K105").Value = 2504.5386
$ws.Range("M105").Value = -757.5385999999999

$ws.Range("H107").Value = 2886.4807
$ws.Range("I107").Value = 2502.3809
$ws.Range("J107").Value = 4499.7
$ws.Range("K107").Value = 2502.3809
$ws.Range("L107").Value = 4499.7
$ws.Range("M107").Value = -582.3809000000001
$ws.Range("N107").Value = -8339.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 853.2
$ws.Range("I16").Value = 691.5
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 691.5
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -404.5
$ws.Range("N16").Value = -2074

$ws.Range("H107").Value = 1157.591
$ws.Range("I107").Value = 752.5333000000001
$ws.Range("J107").Value = 2025.5714
$ws.Range("K107").Value = 752.5333000000001
$ws.Range("L107").Value = 2025.5714
$ws.Range("M107").Value = 1167.4667
$ws.Range("N107").Value = -5865.5714

$ws.Range("H113").Value = 853.2
$ws.Range("I113").Value = 691.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 691.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1478.5
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 6420.1665
$ws.Range("I132").Value = 3642.6667
$ws.Range("K132").Value = 10928.0001
$ws.Range("M132").Value = -8398.000100000001

$ws.Range("H141").Value = 352616.5
$ws.Range("J141").Value = 375129.44
$ws.Range("L141").Value = 375129.44
$ws.Range("N141").Value = -385489.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 153.07143
$ws.Range("I23").Value = 54.6
$ws.Range("J23").Value = 207.77777
$ws.Range("K23").Value = 163.8
$ws.Range("L23").Value = 623.33331
$ws.Range("M23").Value = 71.19999999999999
$ws.Range("N23").Value = -1093.33331

$ws.Range("H61").Value = 167.73334
$ws.Range("I61").Value = 172.91667
$ws.Range("K61").Value = 518.75001
$ws.Range("M61").Value = -303.75001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5291.1333
$ws.Range("I70").Value = 5355.769
$ws.Range("K70").Value = 5355.769
$ws.Range("M70").Value = -5085.769

$ws.Range("H73").Value = 5291.1333
$ws.Range("I73").Value = 5355.769
$ws.Range("K73").Value = 5355.769
$ws.Range("M73").Value = -4419.769

$ws.Range("H80").Value = 26672104
$ws.Range("I80").Value = 40004480
$ws.Range("J80").Value = 7348.5
$ws.Range("K80").Value = 40004480
$ws.Range("L80").Value = 7348.5
$ws.Range("M80").Value = -40003482
$ws.Range("N80").Value = -9344.5

$ws.Range("H83").Value = 26672104
$ws.Range("I83").Value = 40004480
$ws.Range("J83").Value = 7348.5
$ws.Range("K83").Value = 200022400
$ws.Range("L83").Value = 36742.5
$ws.Range("M83").Value = -200017408
$ws.Range("N83").Value = -46726.5

$ws.Range("H113").Value = 3616.3333
$ws.Range("I113").Value = 2243
$ws.Range("J113").Value = 4989.6665
$ws.Range("K113").Value = 2243
$ws.Range("L113").Value = 4989.6665
$ws.Range("M113").Value = -73
$ws.Range("N113").Value = -9329.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1992.3572
$ws.Range("I46").Value = 1100
$ws.Range("J46").Value = 2235.7273
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 2235.7273
$ws.Range("M46").Value = -912
$ws.Range("N46").Value = -2611.7273

$ws.Range("H55").Value = 553.6667
$ws.Range("I55").Value = 476.4
$ws.Range("J55").Value = 940
$ws.Range("K55").Value = 476.4
$ws.Range("L55").Value = 940
$ws.Range("M55").Value = -303.4
$ws.Range("N55").Value = -1286

$ws.Range("H61").Value = 3709.4285
$ws.Range("I61").Value = 3709.4285
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3709.4285
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3507.4285
$ws.Range("N61").ClearContents()

$ws.Range("H82").Value = 2386.182
$ws.Range("I82").Value = 887.25
$ws.Range("J82").Value = 3242.7144
$ws.Range("K82").Value = 887.25
$ws.Range("L82").Value = 3242.7144
$ws.Range("M82").Value = -526.25
$ws.Range("N82").Value = -3964.7144

$ws.Range("H85").Value = 2386.182
$ws.Range("I85").Value = 887.25
$ws.Range("J85").Value = 3242.7144
$ws.Range("K85").Value = 887.25
$ws.Range("L85").Value = 3242.7144
$ws.Range("M85").Value = 360.75
$ws.Range("N85").Value = -5738.7144

$ws.Range("H98").Value = 89354.5
$ws.Range("J98").Value = 89354.5
$ws.Range("L98").Value = 89354.5
$ws.Range("N98").Value = -95344.5

$ws.Range("H100").Value = 66669110
$ws.Range("I100").Value = 142858960
$ws.Range("K100").Value = 142858960
$ws.Range("M100").Value = -142858419

$ws.Range("H113").Value = 3709.4285
$ws.Range("I113").Value = 3709.4285
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3709.4285
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1539.4285
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 10701.583
$ws.Range("I132").Value = 7419.077
$ws.Range("J132").Value = 14580.909
$ws.Range("K132").Value = 22257.231
$ws.Range("L132").Value = 43742.727
$ws.Range("M132").Value = -19727.231
$ws.Range("N132").Value = -48802.727

$ws.Range("H136").Value = 15005414
$ws.Range("I136").Value = 22504774
$ws.Range("K136").Value = 67514322
$ws.Range("M136").Value = -67511772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61332.332
$ws.Range("J46").Value = 61332.332
$ws.Range("L46").Value = 61332.332
$ws.Range("N46").Value = -61794.332

$ws.Range("H100").Value = 461.70834
$ws.Range("I100").Value = 486.7
$ws.Range("J100").Value = 336.75
$ws.Range("K100").Value = 973.4
$ws.Range("L100").Value = 673.5
$ws.Range("M100").Value = -432.4
$ws.Range("N100").Value = -1755.5

$ws.Range("H132").Value = 16600
$ws.Range("I132").Value = 11750.25
$ws.Range("K132").Value = 35250.75
$ws.Range("M132").Value = -32720.75

$ws.Range("H134").Value = 61332.332
$ws.Range("J134").Value = 61332.332
$ws.Range("L134").Value = 183996.996
$ws.Range("N134").Value = -189066.996
